$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)

# Row 42 - first new horizontal-run sample after the previous day's data
$ws.Range("A42").Value = "D20151104T193427"
$ws.Range("B42").Value = 5
$ws.Range("C42").Value = "H"
$ws.Range("D42").Value = 1137.8
$ws.Range("E42").Value = 3.5983000000000001
$ws.Range("F42").Value = 1198
$ws.Range("G42").Value = 318.89999999999998
$ws.Range("H42").Value = 4094

# Row 43 - first vertical sample after horizontal run
$ws.Range("I43").Value = "first vert sample after horz"
$ws.Range("A43").Value = "D20151104T200201"
$ws.Range("J43").Value = "weird shift down in ypos"
$ws.Range("B43").Value = 5
$ws.Range("C43").Value = "V"
$ws.Range("D43").Value = 1267.9000000000001
$ws.Range("E43").Value = 3.4971000000000001
$ws.Range("F43").Value = 1198
$ws.Range("G43").Value = 343.34
$ws.Range("H43").Value = 4434

# Row 44
$ws.Range("A44").Value = "D20151104T202421"
$ws.Range("I44").Value = "pos still weird and low"
$ws.Range("B44").Value = 5
$ws.Range("C44").Value = "V"
$ws.Range("D44").Value = 1251
$ws.Range("E44").Value = 3.5053000000000001
$ws.Range("F44").Value = 1198
$ws.Range("G44").Value = 340.3
$ws.Range("H44").Value = 4385

# Row 45
$ws.Range("A45").Value = "D20151104T204641"
$ws.Range("I45").Value = "pos stilll lower than usual"
$ws.Range("B45").Value = 5
$ws.Range("C45").Value = "V"
$ws.Range("D45").Value = 1185.7
$ws.Range("E45").Value = 3.6105999999999998
$ws.Range("F45").Value = 1198
$ws.Range("G45").Value = 331.4
$ws.Range("H45").Value = 4281

# Row 46 - comment entered before the filename this time
$ws.Range("I46").Value = "pos still low, no junk in sample"
$ws.Range("A46").Value = "D20151104T210900"
$ws.Range("B46").Value = 5
$ws.Range("C46").Value = "V"
$ws.Range("D46").Value = 1222.9000000000001
$ws.Range("E46").Value = 3.5324
$ws.Range("F46").Value = 1198
$ws.Range("G46").Value = 335.3
$ws.Range("H46").Value = 4320

# Row 47
$ws.Range("A47").Value = "D20151104T213356"
$ws.Range("B47").Value = 5
$ws.Range("C47").Value = "V"
$ws.Range("D47").Value = 1165.3
$ws.Range("E47").Value = 3.6246999999999998
$ws.Range("F47").Value = 1198
$ws.Range("G47").Value = 328
$ws.Range("H47").Value = 4224

# Re-point the view at the new bottom of the data and flip the page to
# portrait to match the print setup used for Sheet1.
$win = $excel.ActiveWindow
$ws.Range("E41").Select()

$ps = $ws.PageSetup
$ps.Orientation = 1
